# Error Calculations and Plots
# Update a handful of column-F (error) values and remove two rows
# ("RM 232" and "SC 92") whose data was dropped from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F value corrections ------------------------------------------------
# RM 8  (row 3): previously missing -> now has a computed value
$ws.Range("F3").Value = 17.64
# RM 14 (row 5): previously had a value -> now missing
$ws.Range("F5").Value = ""
# RM 135 (row 21): previously missing -> now has a computed value
$ws.Range("F21").Value = 16.58
# RM 140 (row 23): previously had a value -> now missing
$ws.Range("F23").Value = ""

# --- Remove rows for "RM 232" (row 26) and "SC 92" (row 28) --------------------
# Delete the lower row first so the earlier row index (26) still points at the
# intended row when it is removed afterwards.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# --- SC 193 (now row 32 after the two deletions) gets a computed F value -------
$ws.Range("F32").Value = 17.39
